$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 59
$ws1.Range("F5").Value = 194
$ws1.Range("F7").Value = 9640
$ws1.Range("F11").Value = 2419
$ws1.Range("F12").Value = 162
$ws1.Range("F13").Value = 107
$ws1.Range("F14").Value = 13
$ws1.Range("F16").Value = 275
$ws1.Range("F17").Value = 475
$ws1.Range("F19").Value = 259
$ws1.Range("F20").Value = 1356

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 59
$ws4.Range("F6").Value = 194
$ws4.Range("F8").Value = 9640
$ws4.Range("F12").Value = 2425
$ws4.Range("F13").Value = 162
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 13
$ws4.Range("F17").Value = 275
$ws4.Range("F18").Value = 475
$ws4.Range("F20").Value = 259
$ws4.Range("F21").Value = 1356
